# Update cube metadata: Package4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet / sheet tab to match the new metadata name.
$ws.Name = "Informe-01-010057-A-TC-TM-TP"

# 2. Resize / restructure the columns to the new layout (15 "data" columns
#    instead of 19, with new widths). Widths are expressed in the
#    "character" ColumnWidth units used by the Excel object model; the
#    values below were chosen so the saved OOXML column width attribute
#    lands as close as possible to the target widths.
$ws.Columns.Item(1).ColumnWidth = 26.8567
$ws.Columns.Item(2).ColumnWidth = 43.6667
$ws.Columns.Item(3).ColumnWidth = 17.8267
$ws.Columns.Item(4).ColumnWidth = 54.3667
$ws.Columns.Item(5).ColumnWidth = 33.8067
$ws.Columns.Item(6).ColumnWidth = 35.4767
$ws.Columns.Item(7).ColumnWidth = 46.4467
$ws.Range("H1:J1").EntireColumn.ColumnWidth = 26.8567
$ws.Columns.Item(11).ColumnWidth = 14.6267
$ws.Columns.Item(12).ColumnWidth = 45.6067
$ws.Columns.Item(13).ColumnWidth = 18.3567
$ws.Columns.Item(14).ColumnWidth = 19.2167
$ws.Columns.Item(15).ColumnWidth = 28.6667

# Columns 16-19 drop their old custom widths and go back to a "plain"
# (default-like) width, matching the new 15-column layout.
$ws.Range("P1:S1").EntireColumn.ColumnWidth = 10.6867

# 3. New selection/active cell on the sheet.
$null = $ws.Range("B19").Select()

# 4. Extend the sheet with additional blank rows (7-14) matching the
#    formatting of the existing trailing blank row 6.
for ($r = 7; $r -le 14; $r++) {
  $ws.Rows.Item($r).RowHeight = 12.8
}
